# Update the "Status" column (F) on the BoM tracker sheet.
# Various "To Order" / "To order microSD" / "To order" entries become
# "Ordered" (parts have since been ordered), one item flips to "-" and
# two previously-blank / in-progress rows get filled in as "Have".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F7").Value  = "Have"
$ws.Range("F8").Value  = "Ordered"
$ws.Range("F9").Value  = "Ordered"
$ws.Range("F10").Value = "Ordered"
$ws.Range("F11").Value = "Ordered"
$ws.Range("F12").Value = "Ordered"
$ws.Range("F13").Value = "Ordered"
$ws.Range("F14").Value = "-"
$ws.Range("F15").Value = "Ordered"
$ws.Range("F16").Value = "Ordered"
$ws.Range("F17").Value = "Ordered"
$ws.Range("F19").Value = "Have"
$ws.Range("F21").Value = "Ordered"

# Match the author's final cell selection in the saved workbook.
$ws.Range("K10").Select()
